$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.248.17"
$ws.Range("E2").Value = "  +3.83%  "
$ws.Range("D3").Value = "3.587.69"
$ws.Range("E3").Value = "  +9.15%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'239.47"
$ws.Range("E5").Value = "  +5.67%  "
$ws.Range("D6").Value = "'637.69"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("D7").Value = "'1.49"
$ws.Range("E7").Value = "  +9.71%  "
$ws.Range("D8").Value = "'0.402"
$ws.Range("E8").Value = "  +7.06%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'1.03"
$ws.Range("E10").Value = "  +12.67%  "
$ws.Range("D11").Value = "3.585.82"
$ws.Range("E11").Value = "  +9.16%  "
$ws.Range("D12").Value = "'43.20"
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("E13").Value = "  +5.50%  "
$ws.Range("E14").Value = "  +9.54%  "
$ws.Range("D15").Value = "4.275.29"
$ws.Range("E15").Value = "  +9.99%  "
$ws.Range("D16").Value = "96.181.86"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("E17").Value = "  +6.35%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'13.28"
$ws.Range("E18").Value = "  +25.35%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.569.45"
$ws.Range("E19").Value = "  +8.79%  "
$ws.Range("D20").Value = "'8.03"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "'18.24"
$ws.Range("E21").Value = "  +7.46%  "
$ws.Range("D22").Value = "'0.501"
$ws.Range("E22").Value = "  +15.21%  "
$ws.Range("D23").Value = "'516.86"
$ws.Range("E23").Value = "  +6.67%  "
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'0.0000199"
$ws.Range("E25").Value = "  +12.69%  "
$ws.Range("D26").Value = "'6.68"
$ws.Range("E26").Value = "  +11.81%  "
$ws.Range("D27").Value = "'96.91"
$ws.Range("E27").Value = "  +9.50%  "
$ws.Range("D28").Value = "'12.45"
$ws.Range("E28").Value = "  +7.87%  "
$ws.Range("E29").Value = "  +21.22%  "
$ws.Range("D30").Value = "'0.145"
$ws.Range("E30").Value = "  +6.01%  "
$ws.Range("E31").Value = "  +6.18%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'1.01"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +6.78%  "
$ws.Range("D35").Value = "'30.46"
$ws.Range("E35").Value = "  +9.60%  "
$ws.Range("D36").Value = "'0.567"
$ws.Range("E36").Value = "  +9.19%  "
$ws.Range("D37").Value = "'579.21"
$ws.Range("E37").Value = "  +9.03%  "
$ws.Range("D38").Value = "'7.90"
$ws.Range("E38").Value = "  +9.92%  "
$ws.Range("E39").Value = "  +11.76%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.152"
$ws.Range("E40").Value = "  +4.68%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.930"
$ws.Range("E42").Value = "  +10.21%  "
$ws.Range("D43").Value = "'0.0434"
$ws.Range("E43").Value = "  +7.39%  "
$ws.Range("E44").Value = "  +6.47%  "
$ws.Range("D45").Value = "'23.84"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'5.64"
$ws.Range("E46").Value = "  +7.71%  "
$ws.Range("D47").Value = "'3.54"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  +6.29%  "
$ws.Range("D49").Value = "'53.79"
$ws.Range("E49").Value = "  +4.29%  "
$ws.Range("D50").Value = "'8.17"
$ws.Range("E50").Value = "  +5.69%  "
$ws.Range("E51").Value = "  +5.82%  "
